$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply Q1 advanced estimate data updates to existing rows
$ws.Range("H98").Value = 904863
$ws.Range("M98").Value = 8643085
$ws.Range("M99").Value = 8586172
$ws.Range("M100").Value = 8552545
$ws.Range("M101").Value = 8496269
$ws.Range("M102").Value = 8469971
$ws.Range("M103").Value = 8390310
$ws.Range("M104").Value = 8353735
$ws.Range("M105").Value = 8358612
$ws.Range("M106").Value = 8313955
$ws.Range("M107").Value = 8307250
$ws.Range("M108").Value = 8263899
$ws.Range("M109").Value = 8212233
$ws.Range("H132").Value = 935498
$ws.Range("M132").Value = 14067045
$ws.Range("M133").Value = 14469497
$ws.Range("M134").Value = 14605132
$ws.Range("M135").Value = 14690569
$ws.Range("M136").Value = 14867028
$ws.Range("M137").Value = 14861547
$ws.Range("M138").Value = 14670938
$ws.Range("M139").Value = 14467163
$ws.Range("D140").Value = 823632
$ws.Range("M140").Value = 14131943
$ws.Range("N140").Value = 8935468
$ws.Range("M141").Value = 13894240
$ws.Range("N141").Value = 8781983
$ws.Range("L142").Value = 104.14
$ws.Range("M142").Value = 13750146
$ws.Range("N142").Value = 8776251
$ws.Range("M143").Value = 13622015
$ws.Range("N143").Value = 8836693
$ws.Range("N144").Value = 9042501
$ws.Range("N145").Value = 9395121
$ws.Range("N146").Value = 9715985
$ws.Range("N147").Value = 9916816
$ws.Range("N148").Value = 10176449
$ws.Range("N149").Value = 10402030
$ws.Range("N150").Value = 10586430
$ws.Range("N151").Value = 10834343
$ws.Range("D186").Value = 590445
$ws.Range("H186").Value = 1208607
$ws.Range("M186").Value = 13889397
$ws.Range("N186").Value = 8211250
$ws.Range("M187").Value = 13965630
$ws.Range("N187").Value = 8231241
$ws.Range("M188").Value = 14176863
$ws.Range("N188").Value = 8311027
$ws.Range("M189").Value = 14260171
$ws.Range("N189").Value = 8315387
$ws.Range("M190").Value = 14422498
$ws.Range("N190").Value = 8370968
$ws.Range("D191").Value = 573525
$ws.Range("E191").Value = 207676
$ws.Range("H191").Value = 1072983
$ws.Range("M191").Value = 14527578
$ws.Range("N191").Value = 8368467
$ws.Range("O191").Value = 2551092
$ws.Range("M192").Value = 14532353
$ws.Range("N192").Value = 8290041
$ws.Range("O192").Value = 2547714
$ws.Range("L193").Value = 128.13
$ws.Range("M193").Value = 14760877
$ws.Range("N193").Value = 8365133
$ws.Range("O193").Value = 2572362
$ws.Range("M194").Value = 14843557
$ws.Range("N194").Value = 8388020
$ws.Range("O194").Value = 2570162
$ws.Range("M195").Value = 14888488
$ws.Range("N195").Value = 8325492
$ws.Range("O195").Value = 2571051
$ws.Range("M196").Value = 14967170
$ws.Range("N196").Value = 8319988
$ws.Range("O196").Value = 2563125
$ws.Range("M197").Value = 15089314
$ws.Range("N197").Value = 8341452
$ws.Range("O197").Value = 2566580
$ws.Range("M198").Value = 15141432
$ws.Range("N198").Value = 8330285
$ws.Range("O198").Value = 2567381
$ws.Range("E199").Value = 214625
$ws.Range("M199").Value = 15308976
$ws.Range("N199").Value = 8356533
$ws.Range("O199").Value = 2587930
$ws.Range("M200").Value = 15440585
$ws.Range("N200").Value = 8350552
$ws.Range("O200").Value = 2599983
$ws.Range("E201").Value = 210445
$ws.Range("M201").Value = 15510360
$ws.Range("N201").Value = 8368295
$ws.Range("O201").Value = 2605562
$ws.Range("M202").Value = 15673670
$ws.Range("N202").Value = 8381852
$ws.Range("O202").Value = 2635474
$ws.Range("O203").Value = 2650971
$ws.Range("O204").Value = 2658793
$ws.Range("O205").Value = 2687342
$ws.Range("O206").Value = 2682396
$ws.Range("O207").Value = 2694747
$ws.Range("O208").Value = 2719330
$ws.Range("O209").Value = 2707772
$ws.Range("O210").Value = 2680337
$ws.Range("O211").Value = 2671980
$ws.Range("O212").Value = 2637238
$ws.Range("B605").Value = 17624674
$ws.Range("D605").Value = 12227308
$ws.Range("E605").Value = 206020
$ws.Range("F605").Value = 66075737
$ws.Range("G605").Value = 55690467
$ws.Range("H605").Value = 18332275
$ws.Range("I605").Value = 334.29
$ws.Range("K605").Value = 6.98
$ws.Range("L605").Value = 356.16
$ws.Range("M605").Value = 44510465
$ws.Range("N605").Value = 18301723
$ws.Range("O605").Value = 1812574
$ws.Range("B606").Value = 9279282
$ws.Range("D606").Value = 5756570
$ws.Range("E606").Value = 262760
$ws.Range("F606").Value = 87020528
$ws.Range("G606").Value = 76814427
$ws.Range("H606").Value = 23726957
$ws.Range("I606").Value = 319.44
$ws.Range("J606").Value = 38.01
$ws.Range("K606").Value = 8.390000000000001
$ws.Range("L606").Value = 341.82
$ws.Range("M606").Value = 66231612
$ws.Range("N606").Value = 23707113
$ws.Range("O606").Value = 1924685
$ws.Range("B607").Value = 6646792
$ws.Range("D607").Value = 3159178
$ws.Range("E607").Value = 606117
$ws.Range("F607").Value = 77536350
$ws.Range("G607").Value = 75533295
$ws.Range("H607").Value = 22100209
$ws.Range("I607").Value = 305.89
$ws.Range("J607").Value = 46.67
$ws.Range("K607").Value = 10.15
$ws.Range("L607").Value = 331.51
$ws.Range("M607").Value = 86415971
$ws.Range("N607").Value = 26500891
$ws.Range("O607").Value = 2392117
$ws.Range("B608").Value = 5948513
$ws.Range("D608").Value = 2156113
$ws.Range("E608").Value = 949812
$ws.Range("F608").Value = 68581507
$ws.Range("G608").Value = 64794858
$ws.Range("H608").Value = 18646590
$ws.Range("I608").Value = 303.93
$ws.Range("J608").Value = 62.22
$ws.Range("K608").Value = 11.61
$ws.Range("L608").Value = 325.74
$ws.Range("M608").Value = 102748649
$ws.Range("N608").Value = 28184477
$ws.Range("O608").Value = 3183881
$ws.Range("B609").Value = 3749396
$ws.Range("D609").Value = 1314804
$ws.Range("E609").Value = 1066469
$ws.Range("F609").Value = 63189272
$ws.Range("G609").Value = 56487528
$ws.Range("H609").Value = 16366072
$ws.Range("I609").Value = 304.38
$ws.Range("J609").Value = 80.68000000000001
$ws.Range("K609").Value = 12.97
$ws.Range("L609").Value = 322.12
$ws.Range("M609").Value = 117123685
$ws.Range("N609").Value = 29143106
$ws.Range("O609").Value = 4114407
$ws.Range("B610").Value = 3494289
$ws.Range("D610").Value = 911708
$ws.Range("E610").Value = 3082993
$ws.Range("F610").Value = 49098392
$ws.Range("G610").Value = 45873012
$ws.Range("H610").Value = 13421458
$ws.Range("I610").Value = 308.91
$ws.Range("J610").Value = 109.81
$ws.Range("K610").Value = 14.07
$ws.Range("L610").Value = 320.08
$ws.Range("M610").Value = 128601081
$ws.Range("N610").Value = 29743624
$ws.Range("O610").Value = 7058163
$ws.Range("B611").Value = 3449913
$ws.Range("D611").Value = 761120
$ws.Range("E611").Value = 2433833
$ws.Range("F611").Value = 32048072
$ws.Range("G611").Value = 28413291
$ws.Range("H611").Value = 8356093
$ws.Range("I611").Value = 316
$ws.Range("J611").Value = 51.09
$ws.Range("K611").Value = 14.65
$ws.Range("L611").Value = 319.16
$ws.Range("M611").Value = 135047183
$ws.Range("N611").Value = 30140372
$ws.Range("O611").Value = 9350376
$ws.Range("B612").Value = 3437800
$ws.Range("D612").Value = 741150
$ws.Range("E612").Value = 1329523
$ws.Range("F612").Value = 26429497
$ws.Range("G612").Value = 21910975
$ws.Range("H612").Value = 6687976
$ws.Range("I612").Value = 329.98
$ws.Range("J612").Value = 44.51
$ws.Range("K612").Value = 15.03
$ws.Range("L612").Value = 318.99
$ws.Range("M612").Value = 139909094
$ws.Range("N612").Value = 30514876
$ws.Range("O612").Value = 10551996
$ws.Range("B613").Value = 4226800
$ws.Range("D613").Value = 875107
$ws.Range("E613").Value = 1041335
$ws.Range("F613").Value = 23022902
$ws.Range("G613").Value = 20350461
$ws.Range("H613").Value = 6402239
$ws.Range("I613").Value = 339.43
$ws.Range("J613").Value = 43.15
$ws.Range("K613").Value = 15.3
$ws.Range("L613").Value = 318.84
$ws.Range("M613").Value = 143631987
$ws.Range("N613").Value = 30828546
$ws.Range("O613").Value = 11434420
$ws.Range("B614").Value = 4607942
$ws.Range("D614").Value = 910450
$ws.Range("E614").Value = 821404
$ws.Range("F614").Value = 22604119
$ws.Range("G614").Value = 17582801
$ws.Range("H614").Value = 5641550
$ws.Range("I614").Value = 345.98
$ws.Range("J614").Value = 42.94
$ws.Range("K614").Value = 15.53
$ws.Range("L614").Value = 318.66
$ws.Range("M614").Value = 146326366
$ws.Range("N614").Value = 31002503
$ws.Range("O614").Value = 12103803
$ws.Range("B615").Value = 3556037
$ws.Range("C615").Value = 5029.5
$ws.Range("D615").Value = 765856
$ws.Range("E615").Value = 631956
$ws.Range("F615").Value = 18841932
$ws.Range("G615").Value = 15853945
$ws.Range("H615").Value = 5075617
$ws.Range("I615").Value = 348.94
$ws.Range("K615").Value = 15.66
$ws.Range("L615").Value = 318.52
$ws.Range("M615").Value = 148648937
$ws.Range("N615").Value = 31288929
$ws.Range("O615").Value = 12597112

# Add new row 616 (March 2021 data)
$ws.Range("A616").Value = 44286
$ws.Range("B616").Value = 4195650
$ws.Range("C616").Value = 4460.9
$ws.Range("D616").Value = 864422
$ws.Range("E616").Value = 617296
$ws.Range("F616").Value = 19006560
$ws.Range("G616").Value = 16798299
$ws.Range("H616").Value = 5424557
$ws.Range("I616").Value = 349.6
$ws.Range("J616").Value = 43.87
$ws.Range("K616").Value = 16.3
$ws.Range("L616").Value = 318.22
$ws.Range("M616").Value = 150181595
$ws.Range("N616").Value = 30443786
$ws.Range("O616").Value = 13049518
